$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Files-tab Neo4j query (row 4, column B) was corrected: the `File Type`
# column and the `Breed` column were dropped from the RETURN clause.
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Weimaraner']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# The row shrank once the extra two columns were removed from the query text.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection/scroll moved onto the corrected Files row.
$ws.Range("B4").Select()
